# Recalculate the "Median Value" column (C) relative to the median AFTER
# merging with zip/census tract data, and re-derive the "Tier" column (D)
# from the recalculated values.
#
# The new median divisor is the value that used to live in the row whose
# old ratio equalled the dataset's old median marker (row 16, C16). Dividing
# every score by that value re-bases the whole column so the former marker
# row becomes exactly 1 (i.e. sits exactly at the new median), matching the
# "scores relative to median AFTER merging" re-basing described in the
# commit message.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Find the extent of the data (header in row 1, data starts row 2).
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

# The divisor used to re-base column C: the pre-edit value of C16.
$divisor = $ws.Cells.Item(16, 3).Value2

# Pass 1: divide every Median Value (column C) by the divisor.
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    $old = $cell.Value2
    if ($old -ne $null) {
        $cell.Value2 = $old / $divisor
    }
}

# Pass 2: re-derive the Tier (column D) from the new Median Value.
# Rule (unchanged from the original workbook's own bucketing):
#   - value < 1            -> "Below Median"
#   - value >= 1, ranked in descending order and split into 4 equally
#     sized groups (first group gets any remainder) -> 1st..4th Tier,
#     with 1st Tier being the highest values.
$belowRows = New-Object System.Collections.ArrayList
$aboveObjs = @()

for ($r = 2; $r -le $lastRow; $r++) {
    $v = $ws.Cells.Item($r, 3).Value2
    if ($v -ne $null) {
        if ($v -lt 1) {
            $belowRows.Add($r) | Out-Null
        } else {
            $aboveObjs += [PSCustomObject]@{Row = $r; Val = $v }
        }
    }
}

# Mark "Below Median" rows.
foreach ($r in $belowRows) {
    $ws.Cells.Item($r, 4).Value2 = "Below Median"
}

# Sort the "above median" rows by value descending (highest value = 1st Tier).
$sortedObjs = $aboveObjs | Sort-Object -Property Val -Descending
$sortedRows = @()
foreach ($o in $sortedObjs) {
    $sortedRows += $o.Row
}

$n = $sortedRows.Count
$base = [Math]::Floor($n / 4)
$rem = $n - ($base * 4)

$tierNames = @("1st Tier", "2nd Tier", "3rd Tier", "4th Tier")

$idx = 0
for ($tierIdx = 0; $tierIdx -lt 4; $tierIdx++) {
    $size = $base
    if ($tierIdx -lt $rem) {
        $size = $base + 1
    }
    for ($j = 0; $j -lt $size; $j++) {
        $r = $sortedRows[$idx]
        $ws.Cells.Item($r, 4).Value2 = $tierNames[$tierIdx]
        $idx = $idx + 1
    }
}
